$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44438
$ws.Range("K2").Value = 950
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 975
$ws.Range("P2").Value = 488
$ws.Range("D3").Value = 44616
$ws.Range("J3").Value = 270
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1400
$ws.Range("P3").Value = 700
$ws.Range("D4").Value = 44540
$ws.Range("D6").Value = 44172
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 1400
$ws.Range("P6").Value = 700
$ws.Range("D7").Value = 44392
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 1900
$ws.Range("P7").Value = 950
$ws.Range("D8").Value = 44181
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1100
$ws.Range("N8").Value = '$/atado'
$ws.Range("P8").Value = 1100
$ws.Range("Q8").Value = 1
$ws.Range("D9").Value = 44266
$ws.Range("K9").Value = 1700
$ws.Range("L9").Value = 1800
$ws.Range("M9").Value = 1750
$ws.Range("P9").Value = 875
$ws.Range("D11").Value = 44243
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1250
$ws.Range("P11").Value = 625
$ws.Range("D12").Value = 44525
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1450
$ws.Range("P12").Value = 725
$ws.Range("D13").Value = 44229
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 1800
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 1900
$ws.Range("P13").Value = 950
$ws.Range("D14").Value = 44365
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 1800
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 1900
$ws.Range("P14").Value = 950
$ws.Range("D15").Value = 44202
$ws.Range("D16").Value = 44468
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 950
$ws.Range("N16").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P16").Value = 475
$ws.Range("Q16").Value = 2
$ws.Range("D17").Value = 44390
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2450
$ws.Range("P17").Value = 1225
$ws.Range("D18").Value = 44427
$ws.Range("K18").Value = 1300
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = 1400
$ws.Range("P18").Value = 700
$ws.Range("D19").Value = 44572
$ws.Range("K19").Value = 1400
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = 1450
$ws.Range("P19").Value = 725
$ws.Range("D20").Value = 44385
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 2400
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2450
$ws.Range("P20").Value = 1225
$ws.Range("D21").Value = 44302
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 950
$ws.Range("P21").Value = 475
$ws.Range("D22").Value = 44257
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1400
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = 1450
$ws.Range("P22").Value = 725
$ws.Range("D23").Value = 44447
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = 950
$ws.Range("P23").Value = 475
$ws.Range("D24").Value = 44435
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 900
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = 950
$ws.Range("P24").Value = 475
$ws.Range("D25").Value = 44291
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 1900
$ws.Range("P25").Value = 950
$ws.Range("D26").Value = 44363
$ws.Range("J26").Value = 250
$ws.Range("K26").Value = 2500
$ws.Range("L26").Value = 2800
$ws.Range("M26").Value = 2650
$ws.Range("P26").Value = 1325
$ws.Range("D27").Value = 44403
$ws.Range("K27").Value = 1800
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = 1900
$ws.Range("P27").Value = 950
$ws.Range("D28").Value = 44601
$ws.Range("J28").Value = 270
$ws.Range("K28").Value = 2200
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2350
$ws.Range("P28").Value = 1175
$ws.Range("D29").Value = 44544
$ws.Range("J29").Value = 250
$ws.Range("K29").Value = 900
$ws.Range("L29").Value = 1000
$ws.Range("M29").Value = 950
$ws.Range("P29").Value = 475
